$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = (Get-Date -Year 2024 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0).Date

# Row 67 - new entry (Open / Vincennes)
$ws.Cells.Item(67, 1).Value = 677
$ws.Cells.Item(67, 2).Value = $newDate
$ws.Cells.Item(67, 3).Value = "Open"
$ws.Cells.Item(67, 4).Value = "Vincennes"

# Row 68 - new entry (Open / Vincennes Place)
$ws.Cells.Item(68, 1).Value = 907
$ws.Cells.Item(68, 2).Value = $newDate
$ws.Cells.Item(68, 3).Value = "Open"
$ws.Cells.Item(68, 4).Value = "Vincennes Place"
$ws.Cells.Item(68, 5).Value = 2
$ws.Cells.Item(68, 6).Value = 2

# Row 69 - new entry (Open / Vincennes)
$ws.Cells.Item(69, 1).Value = 907
$ws.Cells.Item(69, 2).Value = $newDate
$ws.Cells.Item(69, 3).Value = "Open"
$ws.Cells.Item(69, 4).Value = "Vincennes"

# Update active selection to reflect the user's last click position
$ws.Range("H66").Select()
